# Bullseye Test Lab Links - replace the old Google Drive hyperlinks in
# column C with new jnfriedman5.github.io links (cell text + hyperlink
# address), in the same order the original author re-added them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bullseye Test Lab Links")

# Drop every existing hyperlink on the sheet so we can re-add them in the
# exact order the workbook's xr:uid / r:id sequence implies.
$ws.Cells.Hyperlinks.Delete()

$ws.Range("C2").Value = "https://jnfriedman5.github.io/BullseyeLabs/?pdf=Target_Essence Flash Poll_27June2025(Sent).pdf"
$ws.Hyperlinks.Add($ws.Range("C2"), "https://jnfriedman5.github.io/BullseyeLabs/?pdf=Target_Essence Flash Poll_27June2025(Sent).pdf")

$ws.Range("C4").Value = "https://jnfriedman5.github.io/BullseyeLabs/?pdf=DEI.pdf"
$ws.Hyperlinks.Add($ws.Range("C4"), "https://jnfriedman5.github.io/BullseyeLabs/?pdf=DEI.pdf")

$ws.Range("C5").Value = "https://jnfriedman5.github.io/BullseyeLabs/?pdf=Taffis.pdf"
$ws.Hyperlinks.Add($ws.Range("C5"), "https://jnfriedman5.github.io/BullseyeLabs/?pdf=Taffis.pdf")

$ws.Range("C3").Value = "https://jnfriedman5.github.io/BullseyeLabs/?pdf=Equality%20Act%202.pdf"
$ws.Hyperlinks.Add($ws.Range("C3"), "https://jnfriedman5.github.io/BullseyeLabs/?pdf=Equality%20Act%202.pdf")

$ws.Range("C7").Value = "https://jnfriedman5.github.io/BullseyeLabs/?pdf=Equality%20Act%201.pdf"
$ws.Hyperlinks.Add($ws.Range("C7"), "https://jnfriedman5.github.io/BullseyeLabs/?pdf=Equality%20Act%201.pdf")

$ws.Range("C6").Value = "https://jnfriedman5.github.io/BullseyeLabs/?pdf=Creative%20Test.pdf"
$ws.Hyperlinks.Add($ws.Range("C6"), "https://jnfriedman5.github.io/BullseyeLabs/?pdf=Creative%20Test.pdf")

# Match the author's final selection state (C4 active).
$ws.Range("C4").Select()
